# Apply the edit described by the diff:
#  - Populate column G (סהכ לזכיין = E*F) and column I (סהכ עמלת רשת = H*E)
#    with formulas for every data row on גיליון1, using a shared formula
#    for the block of rows 4..25 (anchored at row 4) and a standalone
#    formula for row 3 (matching the original author's authoring pattern).
#  - Update the sheet view: scroll position / active selection moved
#    from I28 to I26.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 3: first data row, gets its own (non-shared) formulas ---
$ws.Range("G3").Formula = "=E3*F3"
$ws.Range("I3").Formula = "=H3*E3"

# --- Rows 4-25: shared formula block ---
# Setting the formula across the whole contiguous block first makes Excel
# register a single shared formula (ref="G4:G25" / ref="I4:I25"); the rows
# that have no data in the source sheet (5,7,9,11,12,14,20,22,24) are then
# cleared again so no stray rows get introduced.
$ws.Range("G4:G25").Formula = "=E4*F4"
$ws.Range("I4:I25").Formula = "=H4*E4"

$gapRows = @(5,7,9,11,12,14,20,22,24)
foreach ($r in $gapRows) {
    $ws.Range("G$r").ClearContents()
    $ws.Range("I$r").ClearContents()
}

# --- Sheet view: selection moved from I28 to I26 ---
$ws.Range("I26").Select()
